$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.127.81"
$ws.Range("E2").Value = "  +0.26%  "

$ws.Range("D3").Value = "2.051.58"
$ws.Range("E3").Value = "  -0.47%  "

$ws.Range("E4").Value = "  +0.29%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.669"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.84"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.50%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("E9").Value = "  +1.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0793"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.37%  "

$ws.Range("E11").Value = "  +1.58%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.12"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.97%  "

$ws.Range("D13").Value = "2.350.97"
$ws.Range("E13").Value = "  -0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.833"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.65%  "

$ws.Range("D16").Value = "2.047.14"
$ws.Range("E16").Value = "  -0.69%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.30"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +27.88%  "

$ws.Range("D18").Value = "37.138.72"
$ws.Range("E18").Value = "  +0.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "76.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.16%  "

$ws.Range("D20").Value = "0.0₃0906"
$ws.Range("E20").Value = "  -3.23%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.41"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("E23").Value = "  +0.01%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.61%  "

$ws.Range("E25").Value = "  +10.83%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.11%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.59%  "

$ws.Range("E29").Value = "  +0.84%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +7.30%  "

$ws.Range("E31").Value = "  +5.02%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0631"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.57"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0897"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.08%  "

$ws.Range("E35").Value = "  +0.15%  "

$ws.Range("E36").Value = "  -1.61%  "

$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("E38").Value = "  +3.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.24%  "

$ws.Range("E40").Value = "  +13.62%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.13%  "

$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.65"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.42%  "

$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "97.42"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("E46").Value = "  +3.42%  "

$ws.Range("D47").Value = "1.295.94"
$ws.Range("E47").Value = "  -0.52%  "

$ws.Range("B48").Value = "MXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.21%  "

$ws.Range("B49").Value = "FTXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -10.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.86"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.46%  "

$ws.Range("D51").Value = "2.243.67"
$ws.Range("E51").Value = "  -0.04%  "
